# Update the Summary Info sheet with freshly loaded data before generating reports.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of cell -> new value, derived from the source data refresh.
$updates = @{
    "C2"  = 326
    "M2"  = 38
    "C3"  = 983
    "M3"  = 118
    "C4"  = 637
    "M4"  = 131
    "C5"  = 310
    "M5"  = 54
    "C6"  = 150
    "J6"  = 1
    "M6"  = 28
    "C7"  = 961
    "M7"  = 240
    "C8"  = 936
    "M8"  = 130
    "C9"  = 388
    "M9"  = 90
    "C10" = 496
    "M10" = 85
    "C11" = 453
    "M11" = 105
    "C12" = 401
    "M12" = 74
    "C13" = 111
    "C14" = 154
    "M14" = 16
    "C15" = 743
    "C16" = 886
    "J16" = 412
    "M16" = 115
    "C17" = 592
    "M17" = 120
    "C18" = 754
    "M18" = 145
    "C19" = 626
    "M19" = 109
    "C20" = 510
    "M20" = 80
    "C21" = 1146
    "J21" = 548
    "M21" = 110
    "C22" = 634
    "M22" = 62
    "C23" = 430
    "J23" = 143
    "C24" = 394
    "J24" = 1
    "M24" = 86
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
